# Update solidarity_support_climate_positive data with final rendered values
# (re-run of "prepare & render" produced slightly different numeric results
# for the "All", "Russia", "Saudi Arabia" and "USA" columns on rows 2-6,
# plus a negligible floating point change on C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id 14)
$ws.Range("B2").Value = 0.563498853260559
$ws.Range("K2").Value = 0.449843052433434
$ws.Range("L2").Value = 0.583955154990387
$ws.Range("N2").Value = 0.543153219230958

# Row 3 (id 15)
$ws.Range("B3").Value = 0.548888702456242
$ws.Range("K3").Value = 0.440780635761654
$ws.Range("L3").Value = 0.605453531287603
$ws.Range("N3").Value = 0.51718496787221

# Row 4 (id 16)
$ws.Range("B4").Value = 0.484425269846573
$ws.Range("C4").Value = 0.532496566334214
$ws.Range("K4").Value = 0.323312484205863
$ws.Range("L4").Value = 0.594355780533345
$ws.Range("N4").Value = 0.435731863213516

# Row 5 (id 17)
$ws.Range("B5").Value = 0.472274429131213
$ws.Range("K5").Value = 0.301289618026923
$ws.Range("L5").Value = 0.456838794373566
$ws.Range("N5").Value = 0.460841645436029

# Row 6 (id 18)
$ws.Range("B6").Value = 0.373752935747861
$ws.Range("K6").Value = 0.263792724428881
$ws.Range("L6").Value = 0.34371419848382
$ws.Range("N6").Value = 0.357961152388123
